# chore: update Sheets via scheduled runner
# Refresh scraped market-board figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# in columns H:N for the affected leve rows across the per-crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10120.833
$ws.Range("I33").Value = 2376.4707
$ws.Range("J33").Value = 28928.572
$ws.Range("K33").Value = 2376.4707
$ws.Range("L33").Value = 28928.572
$ws.Range("M33").Value = -2147.4707
$ws.Range("N33").Value = -29386.572
$ws.Range("H111").Value = 1298.1818
$ws.Range("I111").Value = 945
$ws.Range("K111").Value = 2835
$ws.Range("M111").Value = 232
$ws.Range("H132").Value = 1492.1111
$ws.Range("I132").Value = 775.8182
$ws.Range("J132").Value = 3461.9167
$ws.Range("K132").Value = 2327.4546
$ws.Range("L132").Value = 10385.7501
$ws.Range("M132").Value = 202.5454
$ws.Range("N132").Value = -15445.7501
$ws.Range("H133").Value = 38000
$ws.Range("J133").Value = 38000
$ws.Range("L133").Value = 38000
$ws.Range("N133").Value = -48120
$ws.Range("H135").Value = 597.4074000000001
$ws.Range("I135").Value = 297.08334
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 2673.75006
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -138.7500600000003
$ws.Range("N135").Value = -32070
$ws.Range("H138").Value = 1946.6863
$ws.Range("I138").Value = 1002.8108
$ws.Range("J138").Value = 4441.2144
$ws.Range("K138").Value = 3008.4324
$ws.Range("L138").Value = 13323.6432
$ws.Range("M138").Value = 2131.5676
$ws.Range("N138").Value = -23603.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6961.1904
$ws.Range("I32").Value = 3573.3484
$ws.Range("J32").Value = 19383.277
$ws.Range("K32").Value = 3573.3484
$ws.Range("L32").Value = 19383.277
$ws.Range("M32").Value = -3286.3484
$ws.Range("N32").Value = -19957.277
$ws.Range("H45").Value = 1409.3334
$ws.Range("I45").Value = 1342.4
$ws.Range("J45").Value = 1457.1428
$ws.Range("K45").Value = 1342.4
$ws.Range("L45").Value = 1457.1428
$ws.Range("M45").Value = -965.4000000000001
$ws.Range("N45").Value = -2211.1428
$ws.Range("H110").Value = 2603.5
$ws.Range("I110").Value = 930.25
$ws.Range("J110").Value = 5950
$ws.Range("K110").Value = 930.25
$ws.Range("L110").Value = 5950
$ws.Range("M110").Value = 1114.75
$ws.Range("N110").Value = -10040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1509.24
$ws.Range("I134").Value = 1078.7222
$ws.Range("J134").Value = 2616.2856
$ws.Range("K134").Value = 3236.1666
$ws.Range("L134").Value = 7848.8568
$ws.Range("M134").Value = -701.1665999999996
$ws.Range("N134").Value = -12918.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5701.7144
$ws.Range("I16").Value = 4370.3335
$ws.Range("K16").Value = 4370.3335
$ws.Range("M16").Value = -4083.3335
$ws.Range("H51").Value = 8500
$ws.Range("J51").Value = 12000
$ws.Range("L51").Value = 12000
$ws.Range("N51").Value = -13472
$ws.Range("H61").Value = 8500
$ws.Range("J61").Value = 12000
$ws.Range("L61").Value = 12000
$ws.Range("N61").Value = -12696
$ws.Range("H74").Value = 13805.25
$ws.Range("J74").Value = 13805.25
$ws.Range("L74").Value = 13805.25
$ws.Range("N74").Value = -15553.25
$ws.Range("H77").Value = 13805.25
$ws.Range("J77").Value = 13805.25
$ws.Range("L77").Value = 41415.75
$ws.Range("N77").Value = -50151.75
$ws.Range("H107").Value = 1089.9524
$ws.Range("I107").Value = 307.36365
$ws.Range("J107").Value = 1950.8
$ws.Range("K107").Value = 307.36365
$ws.Range("L107").Value = 1950.8
$ws.Range("M107").Value = 1612.63635
$ws.Range("N107").Value = -5790.8
$ws.Range("H113").Value = 5701.7144
$ws.Range("I113").Value = 4370.3335
$ws.Range("K113").Value = 4370.3335
$ws.Range("M113").Value = -2200.3335
$ws.Range("H122").Value = 1553.8286
$ws.Range("I122").Value = 1128.4706
$ws.Range("J122").Value = 1955.5555
$ws.Range("K122").Value = 3385.4118
$ws.Range("L122").Value = 5866.666499999999
$ws.Range("M122").Value = -935.4118000000003
$ws.Range("N122").Value = -10766.6665
$ws.Range("H132").Value = 1158.921
$ws.Range("I132").Value = 945.5833
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2836.7499
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -306.7498999999998
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2941909.2
$ws.Range("I122").Value = 672.55554
$ws.Range("J122").Value = 6250800.5
$ws.Range("K122").Value = 6052.99986
$ws.Range("L122").Value = 56257204.5
$ws.Range("M122").Value = -3602.99986
$ws.Range("N122").Value = -56262104.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("N26").Value = -25560
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -25996
$ws.Range("H107").Value = 924.2222
$ws.Range("I107").Value = 578
$ws.Range("J107").Value = 1616.6666
$ws.Range("K107").Value = 578
$ws.Range("L107").Value = 1616.6666
$ws.Range("M107").Value = 1342
$ws.Range("N107").Value = -5456.6666
$ws.Range("H126").Value = 2172
$ws.Range("I126").Value = 1837.3334
$ws.Range("J126").Value = 2458.8572
$ws.Range("K126").Value = 5512.0002
$ws.Range("L126").Value = 7376.571599999999
$ws.Range("M126").Value = -3042.0002
$ws.Range("N126").Value = -12316.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2718.625
$ws.Range("I7").Value = 1928
$ws.Range("J7").Value = 3193
$ws.Range("K7").Value = 1928
$ws.Range("L7").Value = 3193
$ws.Range("M7").Value = -1816
$ws.Range("N7").Value = -3417
$ws.Range("H55").Value = 324
$ws.Range("I55").Value = 298.66666
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 298.66666
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = -125.66666
$ws.Range("N55").Value = -746
$ws.Range("H93").Value = 10071.357
$ws.Range("I93").Value = 18833.334
$ws.Range("J93").Value = 3499.875
$ws.Range("K93").Value = 18833.334
$ws.Range("L93").Value = 3499.875
$ws.Range("M93").Value = -17585.334
$ws.Range("N93").Value = -5995.875
$ws.Range("H126").Value = 2718.625
$ws.Range("I126").Value = 1928
$ws.Range("J126").Value = 3193
$ws.Range("K126").Value = 5784
$ws.Range("L126").Value = 9579
$ws.Range("M126").Value = -3314
$ws.Range("N126").Value = -14519
$ws.Range("H140").Value = 48333.332
$ws.Range("J140").Value = 48333.332
$ws.Range("L140").Value = 48333.332
$ws.Range("N140").Value = -58693.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 344.4091
$ws.Range("I107").Value = 210.27272
$ws.Range("J107").Value = 478.54544
$ws.Range("K107").Value = 630.81816
$ws.Range("L107").Value = 1435.63632
$ws.Range("M107").Value = 1289.18184
$ws.Range("N107").Value = -5275.63632
